# Daily attendance processing - 2025-12-17 23:52:19
#
# Normalizes the "Recorded By" (column G) values on the active sheet.
# Each cell can contain a comma-separated list of recorder names/emails
# (e.g. "dnasr281@gmail.com, System"). This pass re-orders the names in
# each cell according to a fixed priority so that administrative /
# system accounts consistently sort ahead of regular user accounts,
# while preserving the relative order of any names that are not part
# of the known priority list (stable sort).

function Get-RecorderPriority($name) {
    if ($name.Equals("admin@admin.com")) { return 0 }
    if ($name.Equals("System")) { return 1 }
    if ($name.Equals("system")) { return 2 }
    if ($name.Equals("backup@backdoor.com")) { return 3 }
    if ($name.Equals("dnasr281@gmail.com")) { return 4 }
    return 5
}

function Sort-RecordedBy($raw) {
    $parts = $raw -split "," | ForEach-Object { $_.Trim() }

    $items = @()
    foreach ($name in $parts) {
        $pri = Get-RecorderPriority $name
        $items += [PSCustomObject]@{ Name = $name; Pri = $pri }
    }

    # Sort-Object is a stable sort, so names with equal priority (including
    # any unrecognized ones, which all get the same fallback priority)
    # retain their original relative order.
    $sortedItems = $items | Sort-Object -Property Pri
    $sortedNames = $sortedItems | ForEach-Object { $_.Name }

    return ($sortedNames -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = "Recorded By"
$col = 7

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$r = $firstRow
while ($r -le $lastRow) {
    $cell = $ws.Cells.Item($r, $col)
    $raw = $cell.Value2

    $skip = $false
    if ($null -eq $raw) { $skip = $true }
    if (-not $skip -and -not ($raw -is [string])) { $skip = $true }
    if (-not $skip -and $raw -eq "") { $skip = $true }
    if (-not $skip -and $raw.Equals("Recorded By")) { $skip = $true }
    if (-not $skip -and $raw.IndexOf(",") -lt 0) { $skip = $true }

    if (-not $skip) {
        $newVal = Sort-RecordedBy $raw
        if (-not $newVal.Equals($raw)) {
            $cell.Value = $newVal
        }
    }

    $r = $r + 1
}
